# Sync attendance_reports: reorder "Recorded By" values in column G so the
# entries are listed consistently (swap the first two comma-separated
# names/emails whenever the cell does not already start with "System",
# case-sensitively - PowerShell's -eq/-ceq comparisons are case-insensitive
# in this runtime, so a manual char-code comparison is used instead).

function Test-ExactSystem($s) {
    if ($s.Length -ne 6) { return $false }
    $target = "System"
    for ($i = 0; $i -lt 6; $i++) {
        $c1 = [int][char]$s[$i]
        $c2 = [int][char]$target[$i]
        if ($c1 -ne $c2) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value) { continue }
    if ($value -isnot [string]) { continue }
    if ($value -notlike "*,*") { continue }

    $parts = $value -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    if (-not (Test-ExactSystem $parts[0])) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = [string]::Join(", ", $parts)
    }
}
